$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "pierdoli"
$ws.Range("D1").Value = "smiedzi"

# Row 2
$ws.Range("A2").Value = "'1"
$ws.Range("B2").Value = "debil@gmail.com"
$ws.Range("C2").Value = "no debil no"
$ws.Range("D2").Value = "N/A"

# Row 3
$ws.Range("A3").Value = "'2"
$ws.Range("B3").Value = "idiota@gmail.com"
$ws.Range("C3").Value = "N/A"
$ws.Range("D3").Value = "no idiota no"

# The leading apostrophes above force A2/A3 to stay text (matching the "1"/"2"
# shared-string values in the target), but they also stamp a quote-prefixed
# style on those cells. Put them back on the default "Normal" style so the
# cells don't carry an extra style index.
$ws.Range("A2:A3").Style = "Normal"

# The old E column is gone in the new layout.
$ws.Range("E1").ClearContents()
